# Atualização de bases das ligas, do dia: 14-06-2024 às 20:31
# Rows 298-306 on the "Poland Ekstraklasa" sheet had their match-data rows
# (columns B:AD) reshuffled/updated while the id column (A) stayed sequential.
# We rewrite B:AD for each affected row with the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Poland Ekstraklasa")

$data298 = @(7093820, "Poland Ekstraklasa", 45437.52083333334, "Ruch Chorzow", "Cracovia Krakow", 2, 0, 1, 0, "H", 2.5, 3.4, 2.5, 2.6, 3.6, 2.3, 0, 2.025, 1.825, 3, 2.025, 1.825, 1.6, -1, -1, 1.025, -1, -1, 0.825)
$data299 = @(7083188, "Poland Ekstraklasa", 45437.52083333334, "Legia Warsaw", "Zaglebie Lubin", 2, 1, 2, 0, "H", 1.5, 4, 5.5, 1.6, 4.1, 4.333, -0.75, 1.825, 2.025, 3, 1.875, 1.975, 0.6000000000000001, -1, -1, 0.4125, -0.5, 0, 0)
$data300 = @(7041338, "Poland Ekstraklasa", 45437.52083333334, "Jagiellonia Bialystok", "Warta Poznan", 3, 0, 3, 0, "H", 1.444, 4.75, 5.25, 1.4, 4.75, 5.75, -1.25, 1.9, 1.95, 3, 1.925, 1.925, 0.3999999999999999, -1, -1, 0.8999999999999999, -1, 0, 0)
$data301 = @(7090293, "Poland Ekstraklasa", 45437.52083333334, "Radomiak Radom", "Widzew Lodz", 1, 3, 1, 0, "A", 2.2, 3.1, 3.1, 2.15, 3.2, 3.1, -0.25, 1.925, 1.925, 2.75, 1.9, 1.95, -1, -1, 2.1, -1, 0.925, 0.8999999999999999, -1)
$data302 = @(7083187, "Poland Ekstraklasa", 45437.52083333334, "Lech Poznan", "Korona Kielce", 1, 2, 1, 0, "A", 1.8, 3.8, 3.6, 2.1, 3.7, 2.9, -0.25, 1.9, 1.95, 2.75, 1.925, 1.925, -1, -1, 1.9, -1, 0.95, 0.4625, -0.5)
$data303 = @(7088350, "Poland Ekstraklasa", 45437.52083333334, "Puszcza Niepolomice", "Piast Gliwice", 1, 0, 0, 0, "H", 3, 3.1, 2.3, 2.7, 3, 2.625, 0, 1.975, 1.875, 2.25, 2.025, 1.825, 1.7, -1, -1, 0.9750000000000001, -1, -1, 0.825)
$data304 = @(7083189, "Poland Ekstraklasa", 45437.52083333334, "Pogon Szczecin", "Gornik Zabrze", 1, 0, 0, 0, "H", 1.727, 4, 3.75, 1.55, 4.333, 4.5, -1, 1.925, 1.925, 3.5, 2.025, 1.825, 0.55, -1, -1, 0, 0, -1, 0.825)
$data305 = @(7093821, "Poland Ekstraklasa", 45437.52083333334, "LKS Lodz", "Stal Mielec", 3, 2, 3, 0, "H", 2.5, 3.4, 2.5, 2.2, 3.5, 2.8, -0.25, 2.025, 1.825, 3, 2, 1.85, 1.2, -1, -1, 1.025, -1, 1, -1)
$data306 = @(7074364, "Poland Ekstraklasa", 45437.52083333334, "Rakow Czestochowa", "Slask Wroclaw", 1, 2, 1, 0, "A", 2.5, 3.6, 2.4, 2.15, 3.6, 2.875, -0.25, 1.95, 1.9, 2.5, 1.875, 1.975, -1, -1, 1.875, -1, 0.8999999999999999, 0.875, -1)

$rowsMap = @{
    298 = $data298
    299 = $data299
    300 = $data300
    301 = $data301
    302 = $data302
    303 = $data303
    304 = $data304
    305 = $data305
    306 = $data306
}

foreach ($r in $rowsMap.Keys) {
    $vals = $rowsMap[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, 2 + $i).Value = $vals[$i]
    }
}
